$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1625.3334
$ws.Range("I18").Value = 1750
$ws.Range("J18").Value = 1002
$ws.Range("K18").Value = 1750
$ws.Range("L18").Value = 1002
$ws.Range("M18").Value = -1466
$ws.Range("N18").Value = -1570
$ws.Range("H46").Value = 20782.7
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 20782.7
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 62348.10000000001
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -62586.10000000001
$ws.Range("H51").Value = 5888.8887
$ws.Range("I51").Value = 3550
$ws.Range("K51").Value = 3550
$ws.Range("M51").Value = -3066
$ws.Range("H60").Value = 20782.7
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 20782.7
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 62348.10000000001
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -63316.10000000001
$ws.Range("H62").Value = 2402.1177
$ws.Range("I62").Value = 1997.826
$ws.Range("J62").Value = 3247.4546
$ws.Range("K62").Value = 1997.826
$ws.Range("L62").Value = 3247.4546
$ws.Range("M62").Value = -1373.826
$ws.Range("N62").Value = -4495.4546
$ws.Range("H65").Value = 2402.1177
$ws.Range("I65").Value = 1997.826
$ws.Range("J65").Value = 3247.4546
$ws.Range("K65").Value = 9989.130000000001
$ws.Range("L65").Value = 16237.273
$ws.Range("M65").Value = -6869.130000000001
$ws.Range("N65").Value = -22477.273
$ws.Range("H101").Value = 41667160
$ws.Range("I101").Value = 55555920
$ws.Range("J101").Value = 880
$ws.Range("K101").Value = 166667760
$ws.Range("L101").Value = 2640
$ws.Range("M101").Value = -166666138
$ws.Range("N101").Value = -5884
$ws.Range("H116").Value = 6736.905
$ws.Range("I116").Value = 7619.1177
$ws.Range("J116").Value = 2987.5
$ws.Range("K116").Value = 7619.1177
$ws.Range("L116").Value = 2987.5
$ws.Range("M116").Value = -4177.1177
$ws.Range("N116").Value = -9871.5
$ws.Range("H126").Value = 11994.167
$ws.Range("J126").Value = 11994.167
$ws.Range("L126").Value = 11994.167
$ws.Range("N126").Value = -21874.167
$ws.Range("H135").Value = 819.7143
$ws.Range("I135").Value = 658.5
$ws.Range("K135").Value = 5926.5
$ws.Range("M135").Value = -3391.5
$ws.Range("H137").Value = 1889.7307
$ws.Range("I137").Value = 1475.2632
$ws.Range("J137").Value = 3014.7144
$ws.Range("K137").Value = 4425.7896
$ws.Range("L137").Value = 9044.143199999999
$ws.Range("M137").Value = -1875.7896
$ws.Range("N137").Value = -14144.1432

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 940.8
$ws.Range("I26").Value = 940.8
$ws.Range("K26").Value = 940.8
$ws.Range("M26").Value = -610.8
$ws.Range("H61").Value = 7272.625
$ws.Range("I61").Value = 4581.6665
$ws.Range("J61").Value = 8887.200000000001
$ws.Range("K61").Value = 4581.6665
$ws.Range("L61").Value = 8887.200000000001
$ws.Range("M61").Value = -4369.6665
$ws.Range("N61").Value = -9311.200000000001
$ws.Range("H74").Value = 1226.6471
$ws.Range("I74").Value = 1066.3704
$ws.Range("J74").Value = 1844.8572
$ws.Range("K74").Value = 1066.3704
$ws.Range("L74").Value = 1844.8572
$ws.Range("M74").Value = -192.3704
$ws.Range("N74").Value = -3592.8572
$ws.Range("H77").Value = 1226.6471
$ws.Range("I77").Value = 1066.3704
$ws.Range("J77").Value = 1844.8572
$ws.Range("K77").Value = 5331.852
$ws.Range("L77").Value = 9224.286
$ws.Range("M77").Value = -963.8519999999999
$ws.Range("N77").Value = -17960.286
$ws.Range("H118").Value = 39272.668
$ws.Range("J118").Value = 39272.668
$ws.Range("L118").Value = 39272.668
$ws.Range("N118").Value = -42586.668
$ws.Range("H132").Value = 2087.8223
$ws.Range("I132").Value = 882.8148
$ws.Range("J132").Value = 3895.3333
$ws.Range("K132").Value = 2648.4444
$ws.Range("L132").Value = 11685.9999
$ws.Range("M132").Value = -118.4443999999999
$ws.Range("N132").Value = -16745.9999
$ws.Range("H136").Value = 7272.625
$ws.Range("I136").Value = 4581.6665
$ws.Range("J136").Value = 8887.200000000001
$ws.Range("K136").Value = 13744.9995
$ws.Range("L136").Value = 26661.6
$ws.Range("M136").Value = -11194.9995
$ws.Range("N136").Value = -31761.6

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 80000000
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 80000000
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 80000000
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -80000346
$ws.Range("H82").Value = 29250
$ws.Range("H85").Value = 29250
$ws.Range("H86").Value = 2294
$ws.Range("I86").Value = 2572.6365
$ws.Range("J86").Value = 1783.1666
$ws.Range("K86").Value = 2572.6365
$ws.Range("L86").Value = 1783.1666
$ws.Range("M86").Value = -1449.6365
$ws.Range("N86").Value = -4029.1666
$ws.Range("H89").Value = 2294
$ws.Range("I89").Value = 2572.6365
$ws.Range("J89").Value = 1783.1666
$ws.Range("K89").Value = 12863.1825
$ws.Range("L89").Value = 8915.833000000001
$ws.Range("M89").Value = -7247.182500000001
$ws.Range("N89").Value = -20147.833
$ws.Range("H134").Value = 3441.5
$ws.Range("I134").Value = 3150
$ws.Range("J134").Value = 3483.1428
$ws.Range("K134").Value = 9450
$ws.Range("L134").Value = 10449.4284
$ws.Range("M134").Value = -6915
$ws.Range("N134").Value = -15519.4284

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 745
$ws.Range("I22").Value = 666.6667
$ws.Range("J22").Value = 980
$ws.Range("K22").Value = 666.6667
$ws.Range("L22").Value = 980
$ws.Range("M22").Value = -316.6667
$ws.Range("N22").Value = -1680
$ws.Range("H31").Value = 1686.7878
$ws.Range("I31").Value = 1394.963
$ws.Range("J31").Value = 3000
$ws.Range("K31").Value = 1394.963
$ws.Range("L31").Value = 3000
$ws.Range("M31").Value = -1099.963
$ws.Range("N31").Value = -3590
$ws.Range("H34").Value = 1686.7878
$ws.Range("I34").Value = 1394.963
$ws.Range("J34").Value = 3000
$ws.Range("K34").Value = 1394.963
$ws.Range("L34").Value = 3000
$ws.Range("M34").Value = -1192.963
$ws.Range("N34").Value = -3404
$ws.Range("H99").Value = 5451.423
$ws.Range("I99").Value = 1691.6
$ws.Range("J99").Value = 10578.454
$ws.Range("K99").Value = 1691.6
$ws.Range("L99").Value = 10578.454
$ws.Range("M99").Value = -193.5999999999999
$ws.Range("N99").Value = -13574.454
$ws.Range("H126").Value = 5451.423
$ws.Range("I126").Value = 1691.6
$ws.Range("J126").Value = 10578.454
$ws.Range("K126").Value = 5074.799999999999
$ws.Range("L126").Value = 31735.362
$ws.Range("M126").Value = -2604.799999999999
$ws.Range("N126").Value = -36675.362

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 109.9
$ws.Range("I23").Value = 91.5
$ws.Range("J23").Value = 137.5
$ws.Range("K23").Value = 274.5
$ws.Range("L23").Value = 412.5
$ws.Range("M23").Value = -39.5
$ws.Range("N23").Value = -882.5
$ws.Range("H131").Value = 1065.6578
$ws.Range("J131").Value = 1181.4062
$ws.Range("L131").Value = 3544.2186
$ws.Range("N131").Value = -13624.2186

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3166.6667
$ws.Range("I80").Value = 2000
$ws.Range("J80").Value = 3750
$ws.Range("K80").Value = 2000
$ws.Range("L80").Value = 3750
$ws.Range("M80").Value = -1002
$ws.Range("N80").Value = -5746
$ws.Range("H83").Value = 3166.6667
$ws.Range("I83").Value = 2000
$ws.Range("J83").Value = 3750
$ws.Range("K83").Value = 10000
$ws.Range("L83").Value = 18750
$ws.Range("M83").Value = -5008
$ws.Range("N83").Value = -28734
$ws.Range("H102").Value = 1284.375
$ws.Range("I102").Value = 1332.5454
$ws.Range("J102").Value = 1178.4
$ws.Range("K102").Value = 1332.5454
$ws.Range("L102").Value = 1178.4
$ws.Range("M102").Value = 289.4546
$ws.Range("N102").Value = -4422.4

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7409618
$ws.Range("I122").Value = 12347732
$ws.Range("J122").Value = 2445.8333
$ws.Range("K122").Value = 37043196
$ws.Range("L122").Value = 7337.499899999999
$ws.Range("M122").Value = -37040746
$ws.Range("N122").Value = -12237.4999
$ws.Range("H136").Value = 27781944
$ws.Range("I136").Value = 3495.4285
$ws.Range("J136").Value = 66671770
$ws.Range("K136").Value = 10486.2855
$ws.Range("L136").Value = 200015310
$ws.Range("M136").Value = -7936.2855
$ws.Range("N136").Value = -200020410

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 7750
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 15000
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 15000
$ws.Range("M22").Value = -207
$ws.Range("N22").Value = -15586
$ws.Range("H122").Value = 44890.39
$ws.Range("I122").Value = 54107.05
$ws.Range("J122").Value = 1111.25
$ws.Range("K122").Value = 162321.15
$ws.Range("L122").Value = 3333.75
$ws.Range("M122").Value = -159871.15
$ws.Range("N122").Value = -8233.75
$ws.Range("H126").Value = 57913
$ws.Range("I126").Value = 79221.84
$ws.Range("J126").Value = 2510
$ws.Range("K126").Value = 237665.52
$ws.Range("L126").Value = 7530
$ws.Range("M126").Value = -235195.52
$ws.Range("N126").Value = -12470
